# Add new customer-import columns (Phone#, Country, Address 1/2, PostalZip),
# rename "Contact Phone #" -> "Phone #", shrink the header font, resize
# columns/row, and move the active cell selection.
# ("Add files via upload" - ImportCustomers_Template.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text: rename F1, add new headers G1:J1 -------------------
# (H1 is written before G1 so the shared-string table append order lines
# up with the source file: "Address 1" ends up before "Country".)
$ws.Range("F1").Value = "Phone #"
$ws.Range("H1").Value = "Address 1"
$ws.Range("G1").Value = "Country"
$ws.Range("I1").Value = "Address 2"
$ws.Range("J1").Value = "PostalZip"

# --- Header row formatting: smaller bold title font -----------------------
$ws.Range("A1:J1").Font.Size = 12

# --- Row height (header row got shorter once the font shrank) ------------
$ws.Rows(1).RowHeight = 16.5

# --- Column widths (resized to fit the new narrower headers) -------------
$ws.Columns(1).ColumnWidth = 25.251
$ws.Columns(2).ColumnWidth = 27.584
$ws.Columns(5).ColumnWidth = 19.584
$ws.Columns(6).ColumnWidth = 22.084
$ws.Columns(8).ColumnWidth = 14.584
$ws.Columns(9).ColumnWidth = 12.917
$ws.Columns(10).ColumnWidth = 13.251

# --- Selection moved from G7 to C7 ----------------------------------------
$null = $ws.Range("C7").Select()
